# Add the eurofer low carbon roadmap data:
#  - new shared string "iron ore and steel scrap"
#  - update a handful of cells across the "BF bb", "EAF bb" and "BF EAF bb"
#    sheets to reference it / to the "CO2__emitted" value
#  - move the selected/active tab from "Birat BF" to "BF EAF bb" and
#    update each sheet's remembered selection

$wb = $excel.ActiveWorkbook

# --- "BF bb" sheet (2nd tab) ---------------------------------------------
$wsBFbb = $wb.Worksheets.Item(2)
$wsBFbb.Range("D3").Value = "CO2__emitted"
$wsBFbb.Range("B4").Value = "CO2__emitted"

# --- "EAF bb" sheet (3rd tab) ---------------------------------------------
$wsEAFbb = $wb.Worksheets.Item(3)
$wsEAFbb.Range("D3").Value = "CO2__emitted"
$wsEAFbb.Range("B4").Value = "CO2__emitted"

# --- "BF EAF bb" sheet (4th tab) ------------------------------------------
$wsBFEAFbb = $wb.Worksheets.Item(4)
$wsBFEAFbb.Range("D2").Value = "iron ore and steel scrap"
$wsBFEAFbb.Range("D3").Value = "CO2__emitted"
$wsBFEAFbb.Range("B4").Value = "CO2__emitted"

# --- Update the remembered selections on every sheet ----------------------
$wsBiratBF = $wb.Worksheets.Item(1)
$wsBiratBF.Select()
$wsBiratBF.Range("C18").Select()

$wsBFbb.Select()
$wsBFbb.Range("E5").Select()

$wsEAFbb.Select()
$wsEAFbb.Range("D5").Select()

# Selecting this sheet last makes it the active / visible tab, matching
# the workbook's new activeTab="3" and this sheet's tabSelected="1"
$wsBFEAFbb.Select()
$wsBFEAFbb.Range("D2").Select()
